# JBoss_ACB_Presentation-2015.pptx - minor PPT updates
# (adding standalone-war.war to EAP and OS presos and minor PPT updates)
#
# Slide 5 ("What is EAP? - More things to know"), Content Placeholder:
#   - Paragraph 1: call out that JBoss AS 7 is now "Wildfly"
#   - Paragraph 7: mention manual XML configuration editing as a management option

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 1: "Based on community project JBoss AS 7"
#   -> "Based on community project JBoss AS 7 (Wildfly)"
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Based on community project JBoss AS "
[void]$para1.InsertAfter("7 (")
[void]$para1.InsertAfter("Wildfly")
[void]$para1.InsertAfter(")")

# Paragraph 7: "Managed through its Management Console or Command Line Interface"
#   -> "Managed through its Management Console, Command Line Interface or by manually editing XML configuration"
$para7 = $tr.Paragraphs(7, 1)
$para7.Text = "Managed through its Management "
[void]$para7.InsertAfter("Console, Command ")
[void]$para7.InsertAfter("Line ")
[void]$para7.InsertAfter("Interface or by manually editing XML configuration")
